# Update column F (dSF) values to match freshly re-pulled data.
# (commit: "repull data, push all data, mean calculation")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F5").Value = -1
$ws.Range("F6").Value = 1
$ws.Range("F7").Value = -4
$ws.Range("F8").Value = -8
$ws.Range("F13").Value = -8
$ws.Range("F14").Value = -8
$ws.Range("F17").Value = -6
$ws.Range("F18").Value = -9
$ws.Range("F19").Value = 2
$ws.Range("F21").Value = 0
$ws.Range("F22").Value = -7
$ws.Range("F24").Value = -3
$ws.Range("F25").Value = -3
$ws.Range("F26").Value = -1
$ws.Range("F27").Value = -2
$ws.Range("F28").Value = -1
$ws.Range("F30").Value = -7
$ws.Range("F32").Value = -3
$ws.Range("F34").Value = 4
$ws.Range("F37").Value = -3
